$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.28 = 8598.29 pesos`n✅ 8598.29 pesos = 2.27 = 936.43 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 438
$wsTasas.Range("O10").Value = 3766.05
$wsTasas.Range("N12").Value = 3783
$wsTasas.Range("O12").Value = 412.003
